$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.094.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.07%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.971.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.19%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.010"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'329.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.59%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.04%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4990"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.14%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4219"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'53.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.35%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09256"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +5.24%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.80%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'22.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.92%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.966.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -7.88%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.912"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.27%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.458"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.01%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.010"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.16%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.34%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'91.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.46%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06728"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.59%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'19.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.25%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.008"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.09%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.973"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.29%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'29.115.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.21%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'11.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.35%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.262"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.197.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -6.19%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'20.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.34%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'155.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.99%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.309"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -5.49%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.265"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.36%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'126.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.25%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.051"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.03%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09866"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.35%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.521"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.62%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.827"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.24%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.731"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.26%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02433"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.41%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +3.33%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.06415"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.90%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'9.048"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.6483"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.10%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -1.68%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -2.78%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.11%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.386"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +9.80%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.43%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'13.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.195"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.481"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.66%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00000000322"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.88%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06963"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.28%  "
$ws.Range("E51").Style = "Normal"
